$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, copying the formatting of the other header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Add time_taken values for each data row (F2:F7)
$ws.Range("F2").Value = "2021-10-05 10:51:36.912816"
$ws.Range("F3").Value = "2021-10-05 10:51:36.912827"
$ws.Range("F4").Value = "2021-10-05 10:51:36.912831"
$ws.Range("F5").Value = "2021-10-05 10:51:36.912833"
$ws.Range("F6").Value = "2021-10-05 10:51:36.912836"
$ws.Range("F7").Value = "2021-10-05 10:51:36.912839"
